$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Re-enter the start date (B2) -- this is what causes Excel to mint a new
# cell style (explicit black font color rather than the theme color).
$ws.Range("B2").Value = 44743

# Update the numeric parameters that changed.
$ws.Range("B5").Value = 42
$ws.Range("B6").Value = -5
$ws.Range("B7").Value = 5

# Leave the selection where the user finished editing.
$ws.Range("B4").Select()
